$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to match the repulled data
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = -2
$ws.Range("F5").Value = -2
$ws.Range("F6").Value = -2
$ws.Range("F11").Value = 7
$ws.Range("F12").Value = -4
$ws.Range("F18").Value = -4
$ws.Range("F19").Value = 2
$ws.Range("F20").Value = -4
$ws.Range("F21").Value = -2
$ws.Range("F22").Value = 4
$ws.Range("F23").Value = -3
$ws.Range("F24").Value = 8
$ws.Range("F25").Value = -1
$ws.Range("F27").Value = 1
$ws.Range("F28").Value = 5
$ws.Range("F29").Value = -3
$ws.Range("F30").Value = -3
$ws.Range("F31").Value = -1
$ws.Range("F32").Value = -1
$ws.Range("F33").Value = 2
$ws.Range("F34").Value = 3
$ws.Range("F36").Value = 5
